# Aspects spreadsheet manually merged from aspects branch
# Applies the "TYPE" row/column insertion into the triangular usage matrix on
# the "Category usage" sheet (rows 39-46 / cols K-R), adds the c7 remark to
# the legend (S35/T35), and updates the selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Category usage")

# --- format-copy helper cells (styles reused, none are themselves modified) ---
# style "2"  -> unfilled centered matrix cell (source: C40)
# style "6"  -> quote-prefixed centered "X" cell, no fill (source: P42)
# style "3"  -> bold header cell with fill (source: H39)
# style "5"  -> filled computed/mirror cell (source: B41)
# style "7"  -> filled diagonal cell (source: B40)
# style "11" -> left aligned bold remark text cell (source: T40)

function CopyFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# --- legend: add "c7" remark row ---
$ws.Range("S35").Value = "c7"
$ws.Range("T35").Value = "c1 but for assignments, which are only for readers and creators"

# --- column header row (row 39): shift MULT_I- / MULT_OUT / ASSOC right, insert TYPE ---
$ws.Range("O39").Value = "TYPE"
$ws.Range("P39").Value = "MULT_I-"
$ws.Range("Q39").Value = "MULT_OUT"
CopyFormat "Q39" "R39"
$ws.Range("R39").Value = "ASSOC"

# --- new column R needs the same width as the other matrix columns (L:Q) ---
$ws.Range("R1").ColumnWidth = $ws.Range("Q1").ColumnWidth()

# --- row 40 (REMARK) ---
CopyFormat "P42" "O40"
$ws.Range("O40").Value = "-"
CopyFormat "C40" "R40"
$ws.Range("R40").Value = "-"

# --- row 41 (SORT) ---
CopyFormat "P42" "O41"
$ws.Range("O41").Value = "-"
CopyFormat "C40" "P41"
$ws.Range("P41").Value = "-"
CopyFormat "P42" "Q41"
$ws.Range("Q41").Value = "X"
CopyFormat "C40" "R41"
$ws.Range("R41").Value = "-"

# --- row 42 (LABEL) ---
CopyFormat "P42" "Q42"
$ws.Range("Q42").Value = "X"
CopyFormat "C40" "R42"
$ws.Range("R42").Value = "X"

# --- row 43: was MULT_IN, becomes TYPE ---
$ws.Range("K43").Value = "TYPE"
$ws.Range("Q43").Value = "X"
CopyFormat "C40" "R43"
$ws.Range("R43").Value = "X"

# --- row 44: was MULT_OUT, becomes MULT_IN ---
$ws.Range("K44").Value = "MULT_IN"
CopyFormat "P42" "R44"
$ws.Range("R44").Value = "c1"

# --- row 45: was ASSOC, becomes MULT_OUT ---
$ws.Range("K45").Value = "MULT_OUT"
CopyFormat "P42" "R45"
$ws.Range("R45").Value = "X"

# --- row 46: new ASSOC row ---
$ws.Range("K46").Value = "ASSOC"
CopyFormat "B41" "L46"
$ws.Range("L46").Formula = "=R40"
CopyFormat "B41" "M46"
$ws.Range("M46").Formula = "=R41"
CopyFormat "B41" "N46"
$ws.Range("N46").Formula = "=R42"
CopyFormat "B41" "O46"
$ws.Range("O46").Formula = "=R43"
CopyFormat "B41" "P46"
$ws.Range("P46").Formula = "=R44"
CopyFormat "B41" "Q46"
$ws.Range("Q46").Formula = "=R45"
CopyFormat "B40" "R46"

$excel.CutCopyMode = 0

# --- selection / scroll position ---
$ws.Range("R44").Select()

$wb.Save()
